$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Periodo Mora" labels in column E (rows 16-22) so they run
# chronologically ascending (1905..1911) instead of descending (1911..1905),
# and swap the matching "Valor Mora" amounts in column F (rows 16 and 22)
# so they stay aligned with their (now swapped) period labels.

$ws.Range("E16").Value = "1905"
$ws.Range("E17").Value = "1906"
$ws.Range("E18").Value = "1907"
$ws.Range("E19").Value = "1908"
$ws.Range("E20").Value = "1909"
$ws.Range("E21").Value = "1910"
$ws.Range("E22").Value = "1911"

$ws.Range("F16").Value = 3168
$ws.Range("F22").Value = 6944
